# Natmi following Dr Hou advice:
# Re-run the Nlgn2-Nrxn2 ligand-receptor pair table with a new "ECs" cell
# cluster added alongside the existing "FAPs" and "sCs" clusters, producing
# the full 3x3 sending/target cluster grid (9 data rows instead of 4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flattened row-major table of the new A2:T10 contents (9 rows x 20 cols).
# Columns: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# then the 16 numeric NATMI metric columns (E..T).
$flat = @(
    "ECs", "Nlgn2", "Nrxn2", "ECs", 2, 0.6666666666666666, 1.859580666666667, 5.578742, 0.09135611778511415, 0.09135611778511415, 1, 0.3333333333333333, 0.039512, 0.118536, 0.006183359004302676, 0.006183359004302677, 0.07347575130133334, 0.661281761712, 0.0005648876735047214, 0.0005648876735047215,
    "ECs", "Nlgn2", "Nrxn2", "FAPs", 2, 0.6666666666666666, 1.859580666666667, 5.578742, 0.09135611778511415, 0.09135611778511415, 3, 1, 6.039054, 18.117162, 0.9450708374258476, 0.9450708374258477, 11.230108063356, 101.070972570204, 0.08633800273915221, 0.08633800273915222,
    "ECs", "Nlgn2", "Nrxn2", "sCs", 2, 0.6666666666666666, 1.859580666666667, 5.578742, 0.09135611778511415, 0.09135611778511415, 3, 1, 0.3114883333333333, 0.934465, 0.04874580356984966, 0.04874580356984967, 0.5792376825588889, 5.21313914303, 0.004453227372457224, 0.004453227372457224,
    "FAPs", "Nlgn2", "Nrxn2", "ECs", 3, 1, 10.79119133333334, 32.373574, 0.5301417487077033, 0.5301417487077033, 1, 0.3333333333333333, 0.039512, 0.118536, 0.006183359004302676, 0.006183359004302677, 0.4263815519626667, 3.837433967664001, 0.003278056755428543, 0.003278056755428544,
    "FAPs", "Nlgn2", "Nrxn2", "FAPs", 3, 1, 10.79119133333334, 32.373574, 0.5301417487077033, 0.5301417487077033, 3, 1, 6.039054, 18.117162, 0.9450708374258476, 0.9450708374258477, 65.16858718633202, 586.5172846769881, 0.5010215064055924, 0.5010215064055924,
    "FAPs", "Nlgn2", "Nrxn2", "sCs", 3, 1, 10.79119133333334, 32.373574, 0.5301417487077033, 0.5301417487077033, 3, 1, 0.3114883333333333, 0.934465, 0.04874580356984966, 0.04874580356984967, 3.361330203101112, 30.25197182791, 0.0258421855466823, 0.02584218554668231,
    "sCs", "Nlgn2", "Nrxn2", "ECs", 3, 1, 7.704522333333333, 23.113567, 0.3785021335071827, 0.3785021335071827, 1, 0.3333333333333333, 0.039512, 0.118536, 0.006183359004302676, 0.006183359004302677, 0.3044210864346666, 2.739789777912, 0.002340414575369412, 0.002340414575369412,
    "sCs", "Nlgn2", "Nrxn2", "FAPs", 3, 1, 7.704522333333333, 23.113567, 0.3785021335071827, 0.3785021335071827, 3, 1, 6.039054, 18.117162, 0.9450708374258476, 0.9450708374258477, 46.528026415206, 418.752237736854, 0.3577113282811031, 0.3577113282811032,
    "sCs", "Nlgn2", "Nrxn2", "sCs", 3, 1, 7.704522333333333, 23.113567, 0.3785021335071827, 0.3785021335071827, 3, 1, 0.3114883333333333, 0.934465, 0.04874580356984966, 0.04874580356984967, 2.399868820739444, 21.598819386655, 0.01845039065071014, 0.01845039065071014
)

$numRows = 9
$numCols = 20
for ($i = 0; $i -lt $numRows; $i++) {
    $r = $i + 2
    for ($j = 0; $j -lt $numCols; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $flat[$i * $numCols + $j]
    }
}
